# Update "想去人数" (want-to-go count, column F) figures on the
# "展览" (Worksheets.Item(1)) and "全部类型" (Worksheets.Item(4)) sheets.
# The "演出" and "本地生活" sheets are unaffected by this refresh.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: 展览 ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F5").Value = 457
$ws1.Range("F6").Value = 88
$ws1.Range("F7").Value = 579
$ws1.Range("F8").Value = 81
$ws1.Range("F9").Value = 6904
$ws1.Range("F16").Value = 16372
$ws1.Range("F19").Value = 48
$ws1.Range("F22").Value = 121
$ws1.Range("F23").Value = 11454
$ws1.Range("F25").Value = 1094
$ws1.Range("F28").Value = 391
$ws1.Range("F30").Value = 853

# --- Sheet 4: 全部类型 ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F5").Value = 457
$ws4.Range("F6").Value = 88
$ws4.Range("F7").Value = 579
$ws4.Range("F9").Value = 81
$ws4.Range("F10").Value = 6904
$ws4.Range("F18").Value = 16372
$ws4.Range("F21").Value = 48
$ws4.Range("F24").Value = 121
$ws4.Range("F27").Value = 11454
$ws4.Range("F29").Value = 1094
$ws4.Range("F32").Value = 391
$ws4.Range("F34").Value = 853
